$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..131).
# The update bumps that date by 2 days (2023-09-06 -> 2023-09-08) for all rows.
$ws.Range("C2:C131").Value = 45177
